# Add the results for the 2025-08-23 "Partidos" (matches) session.
# The sheet's data currently ends at row 392; this appends 10 new rows
# (393-402) with the same column layout:
#   fecha | jugador | equipo | posicion | goles | autogoles | arquero |
#   goles_recibidos | tarjetas_amarillas | tarjetas_rojas | asistencias |
#   Penales_Atajados

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$matchDate = "8/23/2025"

# jugador, equipo, posicion, goles, autogoles, arquero, goles_recibidos,
# tarjetas_amarillas, tarjetas_rojas, asistencias, Penales_Atajados
$newRows = @(
    @("Fabian Caicedo",              "Amarillo", "Arquero",       0, 0, $true,  2, 1, 0, 0, 0),
    @("Invitado",                    "Azul",     "Arquero",       0, 0, $true,  4, 0, 0, 0, 0),
    @("Alexander Uribe",             "Amarillo", "Mediocampista", 2, 0, $false, 0, 0, 0, 1, 0),
    @("Andres Tangarife",            "Amarillo", "Delantero",     2, 0, $false, 0, 0, 0, 0, 0),
    @("Sebastian Giraldo",           "Amarillo", "Mediocampista", 0, 0, $false, 0, 0, 0, 2, 0),
    @("Andres Guerrero ",            "Amarillo", "Defensa",       0, 0, $false, 0, 0, 0, 1, 0),
    @("Edwin Casas",                 "Azul",     "Mediocampista", 0, 0, $false, 0, 1, 0, 0, 0),
    @("Arnul David Narvaez",         "Azul",     "Delantero",     1, 0, $false, 0, 0, 0, 1, 0),
    @("Francisco Javier Duran",      "Azul",     "Defensa",       1, 0, $false, 0, 0, 0, 0, 0),
    @("Carlos Fernando Valencia",    "Azul",     "Delantero",     0, 0, $false, 0, 0, 0, 1, 0)
)

$startRow = 393
$row = $startRow
foreach ($data in $newRows) {
    $ws.Cells.Item($row, 1).Value  = $matchDate   # fecha
    $ws.Cells.Item($row, 2).Value  = $data[0]     # jugador
    $ws.Cells.Item($row, 3).Value  = $data[1]     # equipo
    $ws.Cells.Item($row, 4).Value  = $data[2]     # posicion
    $ws.Cells.Item($row, 5).Value  = $data[3]     # goles
    $ws.Cells.Item($row, 6).Value  = $data[4]     # autogoles
    $ws.Cells.Item($row, 7).Value  = $data[5]     # arquero
    $ws.Cells.Item($row, 8).Value  = $data[6]     # goles_recibidos
    $ws.Cells.Item($row, 9).Value  = $data[7]     # tarjetas_amarillas
    $ws.Cells.Item($row, 10).Value = $data[8]     # tarjetas_rojas
    $ws.Cells.Item($row, 11).Value = $data[9]     # asistencias
    $ws.Cells.Item($row, 12).Value = $data[10]    # Penales_Atajados
    $row = $row + 1
}

$lastRow = $row - 1

# Leave the selection on the next empty row, column B, matching where the
# author's cursor ended up after pasting the new data.
$ws.Cells.Item($lastRow + 2, 2).Select() | Out-Null
